# Edit the CVRP results workbook:
#  - add new "Kommentar" entries (J column) for rows that hit the new
#    time_limit=100 handling
#  - add the C5 (Maximal Path Length) value that was missing
#  - insert four new benchmark-instance rows (X-n110-k13, X-n115-k10,
#    X-n120-k6, X-n125-k30) before the X-n129-k18 row, and append a new
#    X-n256-k16 row (instance only, result still pending) at the bottom
#  - fix the error found in the X-n106-k14 row: replace the stray
#    runtime value with the corrected "exceeded" note

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (E-n23-k3): Kommentar re-saved (value unchanged) ---
$ws.Cells.Item(3, 10).Value2 = "time_limit=10, heuristic=10"

# --- Row 4 (E-n30-k3): new Kommentar ---
$ws.Cells.Item(4, 10).Value2 = "time_limit 100 exceeded without paths"

# --- Row 5 (E-n33-k4): add the missing Maximal Path Length + new Kommentar ---
$ws.Cells.Item(5, 3).Value2 = 136
$ws.Cells.Item(5, 10).Value2 = "time_limit=100 exceeded"

# --- Row 7 (E-n76-k14): new Kommentar ---
$ws.Cells.Item(7, 10).Value2 = "time_limit=100 exceeded"

# --- Insert four fresh rows above the X-n129-k18 row (currently row 12) ---
$ws.Range("A12:A15").EntireRow.Insert()

# Row 12: X-n110-k13
$ws.Cells.Item(12, 1).Value2 = "X-n110-k13"
$ws.Cells.Item(12, 2).Value2 = 14971
$ws.Cells.Item(12, 3).Value2 = 16
$ws.Cells.Item(12, 4).Value2 = "13924,6…"

# Row 13: X-n115-k10
$ws.Cells.Item(13, 1).Value2 = "X-n115-k10"
$ws.Cells.Item(13, 2).Value2 = 12747
$ws.Cells.Item(13, 3).Value2 = 171

# Row 14: X-n120-k6
$ws.Cells.Item(14, 1).Value2 = "X-n120-k6"
$ws.Cells.Item(14, 2).Value2 = 13332
$ws.Cells.Item(14, 3).Value2 = 23

# Row 15: X-n125-k30
$ws.Cells.Item(15, 1).Value2 = "X-n125-k30"
$ws.Cells.Item(15, 2).Value2 = 55539
$ws.Cells.Item(15, 3).Value2 = 128

# Row 17: new X-n256-k16 instance, not yet solved
$ws.Cells.Item(17, 1).Value2 = "X-n256-k16"

# Fill in the Maximal Path Length results for the new rows
$ws.Cells.Item(13, 4).Value2 = "11982,7…"
$ws.Cells.Item(14, 4).Value2 = "12308,0..."

# --- Row 11 (X-n106-k14): found the error -> replace D11/E11 with the
#     corrected note, dropping the old runtime value ---
$ws.Cells.Item(11, 4).Value2 = "25186,2..."
$ws.Cells.Item(11, 5).ClearContents()

$ws.Cells.Item(15, 4).Value2 = "54998,2…"
$ws.Cells.Item(16, 4).Value2 = "27830,8…"

$ws.Range("D17").Select()
